$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks (will be recreated below in final row order)
$ws.Cells.Hyperlinks.Delete()

# Clear old data rows 2-9 before rewriting the full, re-sorted list (rows 2-12)
$ws.Range("A2:H9").ClearContents()

# Column H got wider to fit the new skill-summary values (XML width=19 <=> ColumnWidth=18.17 in this engine)
$ws.Columns.Item(8).ColumnWidth = 18.17

# Row 2
$ws.Cells.Item(2, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(2, 2).Value = '産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5450864'
$ws.Cells.Item(2, 7).Value = 383
$ws.Cells.Item(2, 8).Value = '🔥AI,Ai ◆開発'

# Row 3
$ws.Cells.Item(3, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(3, 2).Value = '【自動化】Webサービス更新ツール開発(200アカウント管理)'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5448409'
$ws.Cells.Item(3, 7).Value = 230
$ws.Cells.Item(3, 8).Value = '◆ツール,開発 ◇管理'

# Row 4
$ws.Cells.Item(4, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(4, 2).Value = '【Java/対話システム/心理学実験】協同問題解決プラットフォームの改修開発'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5439921'
$ws.Cells.Item(4, 7).Value = 155
$ws.Cells.Item(4, 8).Value = '★Java ◆開発'

# Row 5
$ws.Cells.Item(5, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(5, 2).Value = '海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5251319'
$ws.Cells.Item(5, 7).Value = 135
$ws.Cells.Item(5, 8).Value = '◆ツール,スクレイピング ◇サイト'

# Row 6
$ws.Cells.Item(6, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(6, 2).Value = '【急募】某新聞社のプロトタイプシステム用チャットボット開発'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5450641'
$ws.Cells.Item(6, 7).Value = 83
$ws.Cells.Item(6, 8).Value = '◆開発'

# Row 7
$ws.Cells.Item(7, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(7, 2).Value = '在宅専業OK│フルスタックエンジニア/開発×データ処理に挑戦!EC運営を支える仕事!'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5450846'
$ws.Cells.Item(7, 7).Value = 75
$ws.Cells.Item(7, 8).Value = '◆開発'

# Row 8
$ws.Cells.Item(8, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(8, 2).Value = '【フルスタックエンジニア募集】新規Webサービス開発'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5450548'
$ws.Cells.Item(8, 7).Value = 75
$ws.Cells.Item(8, 8).Value = '◆開発'

# Row 9
$ws.Cells.Item(9, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(9, 2).Value = '【急募】オンラインガチャサイトに決済機能を導入可能な方'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5450884'
$ws.Cells.Item(9, 7).Value = 33
$ws.Cells.Item(9, 8).Value = '◇サイト'

# Row 10
$ws.Cells.Item(10, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(10, 2).Value = '【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5445466'
$ws.Cells.Item(10, 7).Value = 25

# Row 11
$ws.Cells.Item(11, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(11, 2).Value = '注目 限定公開 PR 限定公開の仕事'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5450323'
$ws.Cells.Item(11, 7).Value = 13

# Row 12
$ws.Cells.Item(12, 1).Value = '2025-12-10 12:39:08'
$ws.Cells.Item(12, 2).Value = '【急募】当社HPのバグ修正をお願いしたいです'
$ws.Cells.Item(12, 3).Value = 'システム開発'
$ws.Cells.Item(12, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = '期限情報なし'
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5450784'
$ws.Cells.Item(12, 7).Value = 10

# Re-create hyperlinks for the URL column, in row order, using the shared "Hyperlink" style
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5450864')
$ws.Cells.Item(2, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5448409')
$ws.Cells.Item(3, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5439921')
$ws.Cells.Item(4, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5251319')
$ws.Cells.Item(5, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5450641')
$ws.Cells.Item(6, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5450846')
$ws.Cells.Item(7, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5450548')
$ws.Cells.Item(8, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), 'https://www.lancers.jp/work/detail/5450884')
$ws.Cells.Item(9, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), 'https://www.lancers.jp/work/detail/5445466')
$ws.Cells.Item(10, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), 'https://www.lancers.jp/work/detail/5450323')
$ws.Cells.Item(11, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), 'https://www.lancers.jp/work/detail/5450784')
$ws.Cells.Item(12, 6).Style = "Hyperlink"
